$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting existing rows 61-173 down to 62-174.
$ws.Rows.Item(61).EntireRow.Insert()

# Populate the new row 61 with data (same Mercado/Categoria context as surrounding rows).
$ws.Range("A61").Value = 9
$ws.Range("B61").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C61").Value = "Metropolitana"
$ws.Range("D61").Value = 44533
$ws.Range("E61").Value = 13
$ws.Range("F61").Value = 100112030
$ws.Range("G61").Value = "Poroto granado"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 52
$ws.Range("K61").Value = 41000
$ws.Range("L61").Value = 43000
$ws.Range("M61").Value = 42000
$ws.Range("N61").Value = "$/malla 25 kilos"
$ws.Range("O61").Value = "Perú"
$ws.Range("P61").Value = 1680
$ws.Range("Q61").Value = 25
$ws.Range("R61").Value = "Hortaliza"
